$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "21.637.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.532.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3957"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3144"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.30"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07136"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -7.62%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.624"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.42"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.570"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.534.82"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001082"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06587"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.07"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.089"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.03%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.79"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.368"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.641.10"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.330"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -8.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.63"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.27"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.841"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.708.65"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.61"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.837"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9360"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -15.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08126"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.478"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -8.61%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.079"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05990"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02191"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.436"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -13.90%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.168"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.34%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.81"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5716"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.97"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.89%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5454"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.160"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "115.68"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.853"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06683"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.14%  "
